# Auto-generated edit script applying cryptos list price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.120.93"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "3.519.88"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'594.85"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'173.81"
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("E9").Value = "  +7.28%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "4.128.24"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'29.11"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").Value = "'0.0000182"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").Value = "67.155.20"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "3.519.17"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").Value = "'397.17"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").Value = "'8.00"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "'73.36"
$ws.Range("D23").Value = "'0.541"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("E25").Value = "  -3.89%  "
$ws.Range("D26").Value = "'10.23"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "'0.181"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "'6.28"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").Value = "'1.45"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").Value = "'24.03"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("D35").Value = "'163.67"
$ws.Range("D36").Value = "'0.895"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("D38").Value = "'6.91"
$ws.Range("E38").Value = "  +3.41%  "
$ws.Range("D39").Value = "'27.79"
$ws.Range("E39").Value = "  +5.53%  "
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "'0.0746"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").Value = "'26.54"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  +3.10%  "
$ws.Range("D44").Value = "2.803.24"
$ws.Range("D45").Value = "'42.82"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").Value = "'339.83"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").Value = "'1.10"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").Value = "'33.50"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "'0.848"
$ws.Range("E51").Value = "  -1.13%  "
